$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for sending clusters that are no longer present
# (ECs = row 2, MuSCs = row 5, Resolving-Mac = row 7).
# Delete from the bottom up so row numbers of earlier rows stay valid.
$ws.Rows("7:7").Delete()
$ws.Rows("5:5").Delete()
$ws.Rows("2:2").Delete()

# Remaining rows are now (in order): FAPs, Inflammatory-Mac, Neutrophils
# -> new rows 2, 3, 4. Update their numeric values to the refreshed TPM data.

# Row 2: FAPs -> Agrp/Mc3r -> FAPs
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.7402036666666666
$ws.Cells.Item(2, 8).Value = 2.220611
$ws.Cells.Item(2, 9).Value = 0.1419214585000345
$ws.Cells.Item(2, 10).Value = 0.1419214585000345
$ws.Cells.Item(2, 17).Value = 0.007170352918999999
$ws.Cells.Item(2, 18).Value = 0.064533176271
$ws.Cells.Item(2, 19).Value = 0.1419214585000345
$ws.Cells.Item(2, 20).Value = 0.1419214585000345

# Row 3: Inflammatory-Mac -> Agrp/Mc3r -> FAPs
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.3561633333333333
$ws.Cells.Item(3, 8).Value = 1.06849
$ws.Cells.Item(3, 9).Value = 0.06828825903893201
$ws.Cells.Item(3, 10).Value = 0.06828825903893201
$ws.Cells.Item(3, 17).Value = 0.00345015421
$ws.Cells.Item(3, 18).Value = 0.03105138789
$ws.Cells.Item(3, 19).Value = 0.06828825903893201
$ws.Cells.Item(3, 20).Value = 0.06828825903893201

# Row 4: Neutrophils -> Agrp/Mc3r -> FAPs
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 4.119219666666667
$ws.Cells.Item(4, 8).Value = 12.357659
$ws.Cells.Item(4, 9).Value = 0.7897902824610334
$ws.Cells.Item(4, 10).Value = 0.7897902824610334
$ws.Cells.Item(4, 17).Value = 0.039902880911
$ws.Cells.Item(4, 18).Value = 0.359125928199
$ws.Cells.Item(4, 19).Value = 0.7897902824610334
$ws.Cells.Item(4, 20).Value = 0.7897902824610334
